$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.325.10"
$ws.Range("E2").Value = "  +0.94%  "
$ws.Range("D3").Value = "1.610.98"
$ws.Range("E3").Value = "  +0.55%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.37%  "
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("E7").Value = "  +0.42%  "
$ws.Range("E8").Value = "  +0.87%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0616"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.36%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.54"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.74%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0814"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.25%  "
$ws.Range("D12").Value = "1.835.13"
$ws.Range("E12").Value = "  +0.54%  "
$ws.Range("D13").Value = "1.610.34"
$ws.Range("E13").Value = "  +0.49%  "
$ws.Range("E14").Value = "  +0.53%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.516"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.16%  "
$ws.Range("D16").Value = "26.294.56"
$ws.Range("E16").Value = "  +0.85%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.33"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.27%  "
$ws.Range("D18").Value = "0.0₃0728"
$ws.Range("E18").Value = "  +1.07%  "
$ws.Range("E19").Value = "  -0.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "202.21"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.21%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.27"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.23%  "
$ws.Range("E23").Value = "  +0.93%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.88"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.97%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.90"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.77%  "
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("E27").Value = "  -0.28%  "
$ws.Range("E28").Value = "  +0.78%  "
$ws.Range("E29").Value = "  +2.44%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0498"
$ws.Range("D30").Style = "Normal"
$ws.Range("E31").Value = "  +0.24%  "
$ws.Range("E32").Value = "  +2.90%  "
$ws.Range("E33").Value = "  -0.19%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.49"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.25%  "
$ws.Range("E35").Value = "  +1.11%  "
$ws.Range("D36").Value = "1.163.20"
$ws.Range("E36").Value = "  +3.39%  "
$ws.Range("E37").Value = "  +2.26%  "
$ws.Range("E38").Value = "  -0.11%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.792"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.83%  "
$ws.Range("E40").Value = "  +0.90%  "
$ws.Range("E41").Value = "  +1.29%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.37"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.23%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.783"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("D44").Value = "1.745.51"
$ws.Range("E44").Value = "  +0.41%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.50"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.50%  "
$ws.Range("E46").Value = "  +1.90%  "
$ws.Range("E47").Value = "  +13.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "53.93"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.06%  "
$ws.Range("E49").Value = "  +0.59%  "
$ws.Range("E50").Value = "  -0.21%  "
$ws.Range("E51").Value = "  -0.18%  "
